$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the diff: (cellRef, newValue)
$updates = @(
    @('D2', '69.919.08'),
    @('E2', '  +1.50%  '),
    @('D3', '3.916.63'),
    @('E3', '  +1.41%  '),
    @('D4', '0.999'),
    @('E4', '  -0.04%  '),
    @('D5', '608.37'),
    @('E5', '  +1.22%  '),
    @('D6', '170.05'),
    @('E6', '  +4.68%  '),
    @('D7', '3.917.24'),
    @('E7', '  +1.44%  '),
    @('D9', '0.537'),
    @('E9', '  +1.05%  '),
    @('D10', '0.169'),
    @('E10', '  +0.78%  '),
    @('E11', '  +1.55%  '),
    @('D12', '0.470'),
    @('E12', '  +2.36%  '),
    @('B13', 'Avalanche'),
    @('C13', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'),
    @('D13', '38.41'),
    @('E13', '  +3.90%  '),
    @('B14', 'ShibaInu'),
    @('C14', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'),
    @('D14', '0.0000256'),
    @('E14', '  +5.05%  '),
    @('D15', '4.576.33'),
    @('E15', '  +1.50%  '),
    @('D16', '3.927.90'),
    @('E16', '  +1.62%  '),
    @('D17', '69.959.40'),
    @('E17', '  +1.33%  '),
    @('D18', '18.82'),
    @('E18', '  +10.20%  '),
    @('D19', '7.64'),
    @('E19', '  +1.17%  '),
    @('D21', '11.22'),
    @('E21', '  -1.37%  '),
    @('D22', '493.64'),
    @('E22', '  +1.78%  '),
    @('D23', '0.748'),
    @('E23', '  +4.08%  '),
    @('D24', '0.0000168'),
    @('E24', '  +3.51%  '),
    @('D25', '85.67'),
    @('E25', '  +2.00%  '),
    @('E26', '  +3.15%  '),
    @('D27', '12.39'),
    @('E27', '  +2.30%  '),
    @('D28', '10.19'),
    @('E28', '  +2.28%  '),
    @('E29', '  +0.11%  '),
    @('E30', '  +1.37%  '),
    @('D31', '4.069.58'),
    @('E31', '  +1.49%  '),
    @('E32', '  +2.90%  '),
    @('D33', '7.87'),
    @('E33', '  -0.62%  '),
    @('D34', '32.18'),
    @('E34', '  -0.42%  '),
    @('D35', '3.881.90'),
    @('E35', '  +1.93%  '),
    @('E36', '  +0.90%  '),
    @('E37', '  +4.29%  '),
    @('D38', '1.05'),
    @('E38', '  +1.42%  '),
    @('E39', '  +1.68%  '),
    @('D40', '3.32'),
    @('E40', '  +11.94%  '),
    @('D41', '0.331'),
    @('E41', '  +3.76%  '),
    @('D42', '1.00'),
    @('E42', '  +0.02%  '),
    @('D43', '2.14'),
    @('E43', '  +7.81%  '),
    @('D44', '439.72'),
    @('E44', '  +0.56%  '),
    @('E45', '  -0.52%  '),
    @('D46', '8.70'),
    @('E46', '  +3.54%  '),
    @('E47', '  -0.01%  '),
    @('D48', '0.0370'),
    @('E48', '  +3.17%  '),
    @('D49', '40.71'),
    @('E49', '  +4.94%  '),
    @('B50', 'Monero'),
    @('C50', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'),
    @('D50', '143.83'),
    @('E50', '  +0.26%  '),
    @('B51', 'FLOKI'),
    @('C51', 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'),
    @('D51', '0.000272'),
    @('E51', '  +19.69%  '),
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $col = $ref -replace "[0-9]+$", ""
    if ($col -eq "D" -or $col -eq "E") {
        $ws.Range($ref).NumberFormat = "@"
    }
    $ws.Range($ref).Value = $val
}

Write-Output "Applied $($updates.Length) cell updates"